# feat: add 2022-Q1 data
#
# - The existing "总计" sheet becomes "2022-Q1" and is repopulated with the
#   per-fund holding breakdown for that quarter (8 columns).
# - A brand new "总计" sheet is appended at the end with the historical
#   per-quarter summary (now holding both 2022-Q1 and 2020-Q4 rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value as literal TEXT (even if it looks like a number),
# without leaving any lingering cell style behind.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 0) Duplicate the existing "总计" sheet and move the copy to the very
#    end of the workbook - this preserves its sheetPr/pageSetup/style
#    info so the new, appended "总计" sheet looks exactly like the
#    original one did. Rename the duplicate out of the way for now.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Copy($null, $lastSheet)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计New"

# ---------------------------------------------------------------------
# 1) Rename the original "总计" sheet to "2022-Q1" and rewrite its
#    contents with the per-fund holdings for the new quarter.
# ---------------------------------------------------------------------
$q1 = $oldTotal
$q1.Name = "2022-Q1"

# Header row - B1:D1 already carry the bold/border style (s=2); reuse it
# for the new E1:H1 headers by copying an existing styled header cell.
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"

$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,5))
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,6))
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,7))
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,8))
$q1.Cells.Item(1,8).Value = "仓位排名"

# Row 2 (A2 already holds the styled index value 0 - leave it as-is).
Set-TextValue $q1.Cells.Item(2,2) "002076"
$q1.Cells.Item(2,3).Value = "浙商中证500指数增强A"
Set-TextValue $q1.Cells.Item(2,4) "14.53"
Set-TextValue $q1.Cells.Item(2,5) "93.68"
Set-TextValue $q1.Cells.Item(2,6) "1.46"
Set-TextValue $q1.Cells.Item(2,7) "0.2121"
$q1.Cells.Item(2,8).Value = 7

# Row 3 (new row) - copy the styled index cell down from A2 first.
$q1.Cells.Item(2,1).Copy($q1.Cells.Item(3,1))
$q1.Cells.Item(3,1).Value = 1
Set-TextValue $q1.Cells.Item(3,2) "007386"
$q1.Cells.Item(3,3).Value = "浙商中证500指数增强C"
Set-TextValue $q1.Cells.Item(3,4) "3.38"
Set-TextValue $q1.Cells.Item(3,5) "93.68"
Set-TextValue $q1.Cells.Item(3,6) "1.46"
Set-TextValue $q1.Cells.Item(3,7) "0.0493"
$q1.Cells.Item(3,8).Value = 7

# ---------------------------------------------------------------------
# 2) Finish the appended "总计" sheet: keep its original 2020-Q4 row,
#    push it to row 3, and insert the 2022-Q1 summary row above it.
# ---------------------------------------------------------------------
$total.Cells.Item(2,1).Copy($total.Cells.Item(3,1))
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2020-Q4"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.01

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.26

$total.Name = "总计"
